# Allow user to activate/disactivate hydropower plants in the list from Excel.
#
# This inserts a new "HPP_active" parameter row right below the header row
# on the "Hydropower plant parameters" sheet, defaulted to active (=1) for
# both existing plants (columns C and D), and updates the named ranges that
# pointed at cells which shifted down by one row as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hydropower plant parameters")

# Insert a new row at row 2 - pushes all the existing parameter rows (and
# their formulas / relative references) down by one, same as using Excel's
# "Insert Sheet Rows" command.
$ws.Rows.Item(2).Insert()

# Populate the new row with the HPP_active parameter.
$ws.Range("A2").Value = "HPP_active"
$ws.Range("B2").Value = "used to include (= 1) or exclude (= 0) plant from current run"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1

# The named ranges referencing this sheet pointed at fixed rows, so they
# need to be shifted down by one to keep referring to the same parameter.
$wb.Names.Item("alpha").RefersTo = "='Hydropower plant parameters'!`$C`$25"
$wb.Names.Item("d_min").RefersTo = "='Hydropower plant parameters'!`$C`$24"
$wb.Names.Item("dP_ramp_turb").RefersTo = "='Hydropower plant parameters'!`$C`$18"
$wb.Names.Item("f_opt").RefersTo = "='Hydropower plant parameters'!`$C`$21"
$wb.Names.Item("f_spill").RefersTo = "='Hydropower plant parameters'!`$C`$22"
$wb.Names.Item("gamma_hydro").RefersTo = "='Hydropower plant parameters'!`$C`$26"
$wb.Names.Item("mu").RefersTo = "='Hydropower plant parameters'!`$C`$27"

# Make "Hydropower plant parameters" the active tab/sheet, with cell E2
# selected (matches the author's saved view state).
$ws.Activate()
$ws.Range("E2").Select()
